$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.229.40"
$ws.Range("E2").Value = "  -0.98%  "
$ws.Range("D3").Value = "2.247.16"
$ws.Range("E3").Value = "  -0.94%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "247.36"
$ws.Range("E5").Value = "  -1.30%  "
$ws.Range("D6").Value = "0.621"
$ws.Range("E6").Value = "  -3.23%  "
$ws.Range("D7").Value = "74.17"
$ws.Range("E7").Value = "  -2.15%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").Value = "0.612"
$ws.Range("E9").Value = "  -4.64%  "
$ws.Range("D10").Value = "42.39"
$ws.Range("E10").Value = "  +5.07%  "
$ws.Range("D11").Value = "0.0939"
$ws.Range("E11").Value = "  -3.40%  "
$ws.Range("E12").Value = "  -2.25%  "
$ws.Range("D13").Value = "0.103"
$ws.Range("E13").Value = "  -3.07%  "
$ws.Range("E14").Value = "  -3.09%  "
$ws.Range("E15").Value = "  -0.73%  "
$ws.Range("D16").Value = "2.240.24"
$ws.Range("E16").Value = "  -1.56%  "
$ws.Range("D17").Value = "42.120.01"
$ws.Range("E17").Value = "  -0.97%  "
$ws.Range("D18").Value = "0.0₃0985"
$ws.Range("E18").Value = "  -0.94%  "
$ws.Range("E19").Value = "  -0.42%  "
$ws.Range("D20").Value = "71.99"
$ws.Range("E20").Value = "  -0.17%  "
$ws.Range("D21").Value = "2.24"
$ws.Range("E21").Value = "  +3.98%  "
$ws.Range("D22").Value = "230.54"
$ws.Range("E22").Value = "  -2.26%  "
$ws.Range("D23").Value = "8.97"
$ws.Range("E23").Value = "  +38.28%  "
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("D25").Value = "11.24"
$ws.Range("E25").Value = "  +0.28%  "
$ws.Range("D26").Value = "3.63"
$ws.Range("E26").Value = "  -5.61%  "
$ws.Range("D27").Value = "2.32"
$ws.Range("E27").Value = "  -2.52%  "
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").Value = "2.16"
$ws.Range("E28").Value = "  -1.71%  "
$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").Value = "169.39"
$ws.Range("E29").Value = "  +1.04%  "
$ws.Range("E30").Value = "  -1.06%  "
$ws.Range("D31").Value = "0.0830"
$ws.Range("E31").Value = "  -2.88%  "
$ws.Range("D32").Value = "0.119"
$ws.Range("E32").Value = "  -4.42%  "
$ws.Range("D33").Value = "30.31"
$ws.Range("E33").Value = "  -5.28%  "
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").Value = "5.26"
$ws.Range("E34").Value = "  +11.03%  "
$ws.Range("B35").Value = "Stellar"
$ws.Range("C35").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D35").Value = "0.125"
$ws.Range("E35").Value = "  -1.71%  "
$ws.Range("D36").Value = "4.53"
$ws.Range("E36").Value = "  -0.10%  "
$ws.Range("E37").Value = "  +0.09%  "
$ws.Range("D38").Value = "13.52"
$ws.Range("E38").Value = "  +0.30%  "
$ws.Range("E39").Value = "  -3.43%  "
$ws.Range("D40").Value = "5.81"
$ws.Range("E40").Value = "  -0.89%  "
$ws.Range("D41").Value = "62.13"
$ws.Range("E41").Value = "  +1.28%  "
$ws.Range("D42").Value = "109.35"
$ws.Range("E42").Value = "  +2.87%  "
$ws.Range("E43").Value = "  -1.74%  "
$ws.Range("E44").Value = "  -2.56%  "
$ws.Range("E45").Value = "  +0.82%  "
$ws.Range("E47").Value = "  -2.65%  "
$ws.Range("E48").Value = "  -0.59%  "
$ws.Range("E49").Value = "  +2.87%  "
$ws.Range("D50").Value = "4.15"
$ws.Range("E50").Value = "  -11.25%  "
$ws.Range("E51").Value = "  -1.17%  "
